$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data values -------------------------------------------------
# Order matters here: it determines the order in which new entries are
# appended to xl/sharedStrings.xml, which must match the target file.

# New "1,2,3"-style sample rows in column AB (rows 7-9), matching the
# style already used by AB6.
$ws.Range("AB7").Value = "null"
$ws.Range("AB8").Value = "{},2,3,4"
$ws.Range("AB9").Value = "vec3,2,2,2"

# New columns AD..AH, each a (name, type) pair stacked in row1/row2,
# matching the existing header (row1) / type (row2) pattern used by the
# rest of the sheet.
$ws.Range("AD1").Value = "c1"
$ws.Range("AD2").Value = "array,array,int"

$ws.Range("AE1").Value = "c2"
$ws.Range("AE2").Value = "list,list,int"

$ws.Range("AF1").Value = "c3"
$ws.Range("AF2").Value = "list,set,int"

$ws.Range("AG1").Value = "c4"
$ws.Range("AG2").Value = "map,int,list,int"

$ws.Range("AH1").Value = "c5"
$ws.Range("AH2").Value = "map,int,map,int,int"

# --- Formatting --------------------------------------------------------
# Reuse the exact formatting already present on sibling cells instead of
# building new styles, by copy/paste-special of formats only.

# AB7:AB9 should look like AB6 (small header-ish font used for the
# "##" sample-value column).
$ws.Range("AB6").Copy()
$ws.Range("AB7:AB9").PasteSpecial(-4122)

# AD1:AH1 should look like AB1 (row-1 field-name style).
$ws.Range("AB1").Copy()
$ws.Range("AD1:AH1").PasteSpecial(-4122)

# AD2:AH2 should look like AB2 (row-2 field-type style).
$ws.Range("AB2").Copy()
$ws.Range("AD2:AH2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Column widths -------------------------------------------------
$ws.Columns.Item(30).ColumnWidth = 11.660714285714286
$ws.Columns.Item(33).ColumnWidth = 12.410714285714286

# --- Selection -------------------------------------------------------
$ws.Range("AE9").Select() | Out-Null
